# Edit script: applies the data changes described in the commit diff
# to the "openTickets" worksheet of wizard_of_oz_experiment_data_open.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Lisa Schmidt / "Analysis for Excel" ticket) ---
# Update the ticketDescriptionHighlighting JSON annotation (column G)
$ws.Range("G2").Value = "[{ start: 143, end: 147, key: `"System`" }​,{ start: 104, end: 128, key: `"Fehlerbeschreibung`" }​,{ start: 67, end: 77, key: `"System`" }]"

# --- Row 3 (Karen Werner / Cisco Softphone ticket) ---
# Update the ticketDescriptionHighlighting JSON annotation (column G)
$ws.Range("G3").Value = "[{ start: 229, end: 297, key: `"Service Anfrage`" },​{ start: 191, end: 192, key: `"System`" }​,{ start: 176, end: 191, key: `"System`" }​,{ start: 129, end: 144, key: `"System`" }]"

# --- Row 4 (Heinz Schubert / Installation Nuance Software ticket) ---
# Shorten the ticketDescription text and rename "Braun-Netzwerk" to "Firmen-Netzwerk" (column F)
$ws.Range("F4").Value = "Installation der Nuance PDF Software kann auf dem Laptop nicht gefunden werden.`n2023-01-08 12:08:29 - Nutzer (Weitere Kommentare)`nMuss ich mit dem VPN verbunden sein? Das Paket wird immer noch nicht im Portal Manager angezeigt. Die Checkpoint-App zur direkten Verbindung mit dem Firmen-Netzwerk ist seit letzter Woche nicht mehr verfügbar."
# Enable wrap text on the (now taller) description cell
$ws.Range("F4").WrapText = $true
# Update the ticketDescriptionHighlighting JSON annotation (column G)
$ws.Range("G4").Value = "[{ start: 130, end: 165, key: `"Auslöser`" }​,{ start: 37, end: 78, key: `"Fehlerbeschreibung`" },​{ start: 24, end: 36, key: `"System`" }]"
# Row 4 becomes taller to fit the wrapped text
$ws.Rows.Item(4).RowHeight = 58

# --- Row 5 (closing "Vielen Dank" experiment-end message) ---
# Collapse the double blank line before the closing sentence into a single line break (column F)
$ws.Range("F5").Value = "Bitte füllen Sie im Anschluss nun folgende <a href=`"https://forms.gle/QxKwBcwiMtDbRGhd6`" >`"Umfrage`"</a> aus. Im Anschluss würden wir Sie bitten an einem kurzen Interview teilzunehmen. Falls Sie noch keinen Termin haben sollten, kommen Sie hier zur <a href=`"https://calendly.com/philipp-reinhard-1/interview_experiment`" >`"Terminbuchung`"</a> .`nVielen Dank, dass Sie unsere Forschung im Bereich künstlicher Intelligenz unterstützen!"

# Renumber the trailing systemId helper columns (H:M) on row 5 from 18-23 down to 0-5
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 4
$ws.Range("M5").Value = 5

# --- Sheet view: update the active selection to F4 and scroll one column to the right ---
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F4").Select()
